{"js": "// 1. Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. Split the mailing-address paragraph \"969 Story Road, San Jose CA 95122\"\n//    (the standalone one near the top of the letter, not the identical text\n//    that also appears inside the info table further down) into two\n//    separate paragraphs: \"969 Story Road\" and \"San Jose, CA 95122\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst addressCandidates = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"969 Story Road, San Jose CA 95122\") {\n    addressCandidates.push(paragraphs.items[i]);\n  }\n}\nfor (const p of addressCandidates) {\n  p.parentTableOrNullObject.load(\"isNullObject\");\n}\nawait context.sync();\n\nlet addressPara = null;\nfor (const p of addressCandidates) {\n  if (p.parentTableOrNullObject.isNullObject) {\n    addressPara = p;\n    break;\n  }\n}\n\nif (addressPara) {\n  // Insert the new second line right after the existing paragraph, matching\n  // its formatting (same paragraph/run style is inherited automatically).\n  addressPara.insertParagraph(\"San Jose, CA 95122\", Word.InsertLocation.after);\n  await context.sync();\n\n  // Trim the original paragraph down to just the street address.\n  addressPara.insertText(\"969 Story Road\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Remove the blank \"No Spacing\" paragraph that immediately follows the\n//    \"...Board of Directors\" signature line.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet boardParaIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"Board of Directors\") !== -1) {\n    boardParaIndex = i;\n    break;\n  }\n}\n\nif (boardParaIndex !== -1 && boardParaIndex + 1 < paragraphs2.items.length) {\n  const blankPara = paragraphs2.items[boardParaIndex + 1];\n  blankPara.load(\"text\");\n  await context.sync();\n  if (blankPara.text === \"\") {\n    blankPara.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$find = $d.Content.Find\n$find.Execute(\"September 19, 2025\", $false, $false, $false, $false, $false, $true, 1, $false, \"September 21, 2025\", 2)\n\n# 2. Split the mailing-address paragraph \"969 Story Road, San Jose CA 95122\"\n#    (the standalone one, not the copy inside the info table) into two\n#    paragraphs: \"969 Story Road\" and \"San Jose, CA 95122\".\n$addressPara = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*969 Story Road, San Jose CA 95122*\" -and $p.Range.Information(12) -eq $false) {\n    $addressPara = $p\n    break\n  }\n}\nif ($addressPara -ne $null) {\n  $addressPara.Range.Text = \"969 Story Road`rSan Jose, CA 95122\"\n}\n\n# 3. Remove the blank \"No Spacing\" paragraph that immediately follows the\n#    \"...Board of Directors\" signature line.\n$boardPara = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*Board of Directors*\") {\n    $boardPara = $p\n    break\n  }\n}\nif ($boardPara -ne $null) {\n  $blankPara = $boardPara.Next()\n  if ($blankPara.Range.Text.Trim() -eq \"\") {\n    $blankPara.Range.Delete()\n  }\n}\n"}
